# Fix the farm id (id_azienda) value for the row that currently holds 5.
# The "Comune"/"Localita" for that row is Casalfiumanese, whose id_azienda
# should be 10, not 5.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")
$ws.Range("A6").Value = 10
